# Generate Report for Archive
# The localization status moved on from "Ready for handoff" to "In Translation"
# for the zh-cn / de-de targets, so update every cell that carries that status
# (the Overview roll-up columns plus each language sheet's own Status column),
# then let Excel re-fit the Status columns now that the text is shorter.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# After the text shrinks, re-fit the Status column(s) to the new, narrower
# content (mirrors what Excel's own "AutoFit Column Width" does when the
# report generator regenerates this sheet).
$fitWidth = 12.5

# --- Overview sheet: one status column per locale (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $fitWidth
$wsOverview.Columns.Item(6).ColumnWidth = $fitWidth

# --- zh-cn sheet: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $fitWidth

# --- de-de sheet: Status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $fitWidth
